$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the values in row 3 and row 4 for columns D, J, K, L, M, O, P
$cols = @("D", "J", "K", "L", "M", "O", "P")

foreach ($col in $cols) {
    $cell3 = $ws.Range($col + "3")
    $cell4 = $ws.Range($col + "4")
    $tmp = $cell3.Value2
    $cell3.Value2 = $cell4.Value2
    $cell4.Value2 = $tmp
}
